{"js": "// Update the date paragraph and the five rows of division problems in the\n// single table. The table has 5 data rows (1 per \"week\"), each with 5 cells,\n// interleaved with empty spacer rows. We address cells directly via\n// table.getCell(row, col) using the row indices that actually hold text.\n\nconst body = context.document.body;\n\n// Load paragraphs (first paragraph holds the date) and tables.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// 1) Date paragraph: \"2024-04-07 Sunday\" -> \"2024-04-08 Monday\"\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text.trim() === \"2024-04-07 Sunday\") {\n  dateParagraph.insertText(\"2024-04-08 Monday\", \"Replace\");\n}\n\n// 2) Table cells: old -> new values, given row-major in document order.\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst replacements = [\n  [\"36\u00f74=9, 0\", \"58\u00f72=29, 0\"],\n  [\"21\u00f77=3, 0\", \"59\u00f79=6, 5\"],\n  [\"16\u00f75=3, 1\", \"29\u00f78=3, 5\"],\n  [\"13\u00f76=2, 1\", \"88\u00f75=17, 3\"],\n  [\"23\u00f72=11, 1\", \"43\u00f73=14, 1\"],\n  [\"61\u00f79=6, 7\", \"46\u00f78=5, 6\"],\n  [\"30\u00f74=7, 2\", \"33\u00f79=3, 6\"],\n  [\"62\u00f74=15, 2\", \"59\u00f72=29, 1\"],\n  [\"12\u00f73=4, 0\", \"58\u00f73=19, 1\"],\n  [\"14\u00f73=4, 2\", \"71\u00f78=8, 7\"],\n  [\"47\u00f73=15, 2\", \"76\u00f76=12, 4\"],\n  [\"58\u00f72=29, 0\", \"39\u00f76=6, 3\"],\n  [\"14\u00f76=2, 2\", \"11\u00f76=1, 5\"],\n  [\"27\u00f72=13, 1\", \"69\u00f75=13, 4\"],\n  [\"86\u00f76=14, 2\", \"71\u00f75=14, 1\"],\n  [\"91\u00f78=11, 3\", \"48\u00f74=12, 0\"],\n  [\"14\u00f79=1, 5\", \"25\u00f77=3, 4\"],\n  [\"73\u00f73=24, 1\", \"14\u00f72=7, 0\"],\n  [\"67\u00f77=9, 4\", \"45\u00f75=9, 0\"],\n  [\"94\u00f77=13, 3\", \"51\u00f74=12, 3\"],\n  [\"89\u00f74=22, 1\", \"91\u00f75=18, 1\"],\n  [\"61\u00f79=6, 7\", \"80\u00f75=16, 0\"],\n  [\"53\u00f72=26, 1\", \"32\u00f79=3, 5\"],\n  [\"51\u00f76=8, 3\", \"42\u00f74=10, 2\"],\n  [\"60\u00f74=15, 0\", \"98\u00f73=32, 2\"],\n];\n\n// Find the rows that actually contain text (5 cells each) and walk them in\n// document order, applying the replacements in sequence. This sidesteps the\n// duplicate-text problem (\"61\u00f79=6, 7\" appears twice with different targets).\nlet idx = 0;\nfor (const row of table.rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  const cells = row.cells.items;\n  // Load each cell's first paragraph (holds the run with the number text)\n  // up front, so the formatted run (font/size) survives the replace -\n  // replacing text on the paragraph (or its range), rather than the whole\n  // cell body, keeps the existing <w:rPr> instead of resetting to defaults.\n  for (const cell of cells) {\n    cell.body.paragraphs.load(\"items\");\n  }\n  await context.sync();\n\n  const firstParagraphs = cells.map((c) => c.body.paragraphs.items[0]);\n  for (const p of firstParagraphs) {\n    p.load(\"text\");\n  }\n  await context.sync();\n\n  const hasText = firstParagraphs.some((p) => p.text.trim().length > 0);\n  if (!hasText) continue;\n\n  for (const p of firstParagraphs) {\n    const current = p.text.trim();\n    const [expectedOld, newVal] = replacements[idx];\n    if (current === expectedOld) {\n      p.insertText(newVal, \"Replace\");\n    }\n    idx++;\n  }\n  await context.sync();\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the five rows of division problems in the\n# single table. The table has 5 data rows (rows 1,5,9,13,17 of the 20-row\n# table; the rows in between are blank spacer rows), each with 5 columns.\n\n$d = $word.ActiveDocument\n\n# 1) Date paragraph: \"2024-04-07 Sunday\" -> \"2024-04-08 Monday\"\n$dateParagraph = $d.Paragraphs.Item(1)\nif ($dateParagraph.Range.Text.TrimEnd(\"`r\") -eq \"2024-04-07 Sunday\") {\n    $dateParagraph.Range.Text = \"2024-04-08 Monday\"\n}\n\n# 2) Table cells: old -> new values, addressed by (row, col) in the single\n#    table (1-based), walked top-to-bottom / left-to-right, matching the\n#    document order of the diff.\n$t = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$replacements = @(\n    @(\"36\u00f74=9, 0\", \"58\u00f72=29, 0\"),\n    @(\"21\u00f77=3, 0\", \"59\u00f79=6, 5\"),\n    @(\"16\u00f75=3, 1\", \"29\u00f78=3, 5\"),\n    @(\"13\u00f76=2, 1\", \"88\u00f75=17, 3\"),\n    @(\"23\u00f72=11, 1\", \"43\u00f73=14, 1\"),\n    @(\"61\u00f79=6, 7\", \"46\u00f78=5, 6\"),\n    @(\"30\u00f74=7, 2\", \"33\u00f79=3, 6\"),\n    @(\"62\u00f74=15, 2\", \"59\u00f72=29, 1\"),\n    @(\"12\u00f73=4, 0\", \"58\u00f73=19, 1\"),\n    @(\"14\u00f73=4, 2\", \"71\u00f78=8, 7\"),\n    @(\"47\u00f73=15, 2\", \"76\u00f76=12, 4\"),\n    @(\"58\u00f72=29, 0\", \"39\u00f76=6, 3\"),\n    @(\"14\u00f76=2, 2\", \"11\u00f76=1, 5\"),\n    @(\"27\u00f72=13, 1\", \"69\u00f75=13, 4\"),\n    @(\"86\u00f76=14, 2\", \"71\u00f75=14, 1\"),\n    @(\"91\u00f78=11, 3\", \"48\u00f74=12, 0\"),\n    @(\"14\u00f79=1, 5\", \"25\u00f77=3, 4\"),\n    @(\"73\u00f73=24, 1\", \"14\u00f72=7, 0\"),\n    @(\"67\u00f77=9, 4\", \"45\u00f75=9, 0\"),\n    @(\"94\u00f77=13, 3\", \"51\u00f74=12, 3\"),\n    @(\"89\u00f74=22, 1\", \"91\u00f75=18, 1\"),\n    @(\"61\u00f79=6, 7\", \"80\u00f75=16, 0\"),\n    @(\"53\u00f72=26, 1\", \"32\u00f79=3, 5\"),\n    @(\"51\u00f76=8, 3\", \"42\u00f74=10, 2\"),\n    @(\"60\u00f74=15, 0\", \"98\u00f73=32, 2\")\n)\n\n$idx = 0\nforeach ($row in $dataRows) {\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $t.Cell($row, $col)\n        $expectedOld = $replacements[$idx][0]\n        $newVal = $replacements[$idx][1]\n        $current = $cell.Range.Text.TrimEnd(\"`r\", [char]7)\n        if ($current -eq $expectedOld) {\n            $cell.Range.Text = $newVal\n        }\n        $idx++\n    }\n}\n"}
